# Daily attendance processing - 2026-01-02 11:53:56
# Rotate the "Recorded By" (column G) contributor list for the rows that
# were re-processed today: the most-recently-added recorder (the last
# name in the comma-separated list) is promoted to the front of the list.
#
# Example:  "System, dnasr281@gmail.com"  ->  "dnasr281@gmail.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows in the "Session Analysis Results" sheet whose "Recorded By" column
# (G) needs its contributor order rotated during today's run.
$rowsToProcess = @(
    2, 3, 6, 10, 11, 12, 13, 14, 15, 17, 18, 19, 20, 21, 22, 24, 26, 28, 29,
    32, 36, 37, 38, 39, 40, 41, 43, 44, 45, 46, 47, 48, 50, 52, 54, 55, 58,
    62, 63, 64, 65, 66, 67, 69, 70, 71, 72, 73, 74, 76, 78, 83, 84, 85, 86,
    87, 90, 92, 93, 94, 96, 99, 101, 109, 110, 111, 112, 113, 116, 118, 119,
    120, 122, 125, 127, 135, 136, 137, 138, 139, 142, 144, 145, 146, 148,
    151, 153
)

foreach ($row in $rowsToProcess) {
    $cell = $ws.Cells.Item($row, 7)   # column G = "Recorded By"
    $current = [string]$cell.Value2

    if ([string]::IsNullOrEmpty($current)) {
        continue
    }

    $parts = $current -split ',\s*'
    if ($parts.Count -le 1) {
        continue
    }

    # Move the last recorder to the front of the list.
    $lastIndex = $parts.Count - 1
    $rotated = @($parts[$lastIndex]) + $parts[0..($lastIndex - 1)]

    $cell.Value = [string]::Join(', ', $rotated)
}
